# Updated cryptos list on Mon Jul 31 22:36:50 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.220.08"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "1.855.68"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.51"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6966"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07767"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3072"
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.80"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07814"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "1.855.69"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.106"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.10"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6877"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.519"
$ws.Range("E16").Value = "  +2.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008435"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "29.219.17"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "248.86"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("D20").Value = "2.107.34"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.522"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1493"
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.23"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.870"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.50"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  +4.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.244"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.206"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("E33").Value = "  -0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7595"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.842"
$ws.Range("E36").Value = "  -1.91%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "1.221.26"
$ws.Range("E39").Value = "  -1.78%  "
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9005"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "109.96"
$ws.Range("E42").Value = "  -0.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9991"
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.510"
$ws.Range("E44").Value = "  -11.61%  "
$ws.Range("D45").Value = "2.007.32"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "65.45"
$ws.Range("E47").Value = "  -8.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.575"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5177"
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.753"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.043"
$ws.Range("E51").Value = "  +0.80%  "
